$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern: for the Price column (D), values look numeric (e.g. "210.91")
# and Excel/COM would silently coerce them to a Number type, changing the
# underlying cell type from Text to Number (and losing formatting like
# trailing zeros, e.g. "23.60" -> 23.6). We force literal text by switching
# the cell to a Text number format before the assignment, then restore the
# original cell style afterwards so no visible style/format changes leak in.

$savedStyle = $ws.Cells.Item(2, 4).Style
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '28.351.16'
$ws.Cells.Item(2, 4).Style = $savedStyle
$ws.Cells.Item(2, 5).Value = '  -0.44%  '

$savedStyle = $ws.Cells.Item(3, 4).Style
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '1.566.20'
$ws.Cells.Item(3, 4).Style = $savedStyle
$ws.Cells.Item(3, 5).Value = '  +0.07%  '

$ws.Cells.Item(4, 5).Value = '  -0.15%  '

$savedStyle = $ws.Cells.Item(5, 4).Style
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '210.91'
$ws.Cells.Item(5, 4).Style = $savedStyle
$ws.Cells.Item(5, 5).Value = '  -0.42%  '

$savedStyle = $ws.Cells.Item(6, 4).Style
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.489'
$ws.Cells.Item(6, 4).Style = $savedStyle
$ws.Cells.Item(6, 5).Value = '  -0.60%  '

$ws.Cells.Item(7, 5).Value = '  -0.14%  '

$savedStyle = $ws.Cells.Item(8, 4).Style
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '44.47'
$ws.Cells.Item(8, 4).Style = $savedStyle
$ws.Cells.Item(8, 5).Value = '  -3.69%  '

$savedStyle = $ws.Cells.Item(9, 4).Style
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '23.60'
$ws.Cells.Item(9, 4).Style = $savedStyle
$ws.Cells.Item(9, 5).Value = '  -1.67%  '

$ws.Cells.Item(10, 5).Value = '  -1.29%  '

$ws.Cells.Item(11, 5).Value = '  -0.90%  '

$ws.Cells.Item(12, 5).Value = '  +0.81%  '

$savedStyle = $ws.Cells.Item(13, 4).Style
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '1.788.93'
$ws.Cells.Item(13, 4).Style = $savedStyle
$ws.Cells.Item(13, 5).Value = '  -0.02%  '

$savedStyle = $ws.Cells.Item(14, 4).Style
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '1.564.96'
$ws.Cells.Item(14, 4).Style = $savedStyle
$ws.Cells.Item(14, 5).Value = '  -0.07%  '

$ws.Cells.Item(15, 2).Value = 'WrappedBTC'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$savedStyle = $ws.Cells.Item(15, 4).Style
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '28.377.52'
$ws.Cells.Item(15, 4).Style = $savedStyle
$ws.Cells.Item(15, 5).Value = '  -0.46%  '

$ws.Cells.Item(16, 2).Value = 'Polkadot'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$savedStyle = $ws.Cells.Item(16, 4).Style
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '3.66'
$ws.Cells.Item(16, 4).Style = $savedStyle
$ws.Cells.Item(16, 5).Value = '  -0.33%  '

$ws.Cells.Item(17, 5).Value = '  -1.17%  '

$savedStyle = $ws.Cells.Item(18, 4).Style
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '60.95'
$ws.Cells.Item(18, 4).Style = $savedStyle
$ws.Cells.Item(18, 5).Value = '  -1.99%  '

$savedStyle = $ws.Cells.Item(19, 4).Style
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '228.15'
$ws.Cells.Item(19, 4).Style = $savedStyle

$ws.Cells.Item(20, 5).Value = '  +0.57%  '

$ws.Cells.Item(21, 5).Value = '  -1.88%  '

$ws.Cells.Item(22, 5).Value = '  -0.03%  '

$ws.Cells.Item(23, 5).Value = '  +1.66%  '

$savedStyle = $ws.Cells.Item(24, 4).Style
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '8.94'
$ws.Cells.Item(24, 4).Style = $savedStyle
$ws.Cells.Item(24, 5).Value = '  -2.12%  '

$ws.Cells.Item(25, 5).Value = '  -1.79%  '

$savedStyle = $ws.Cells.Item(26, 4).Style
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '150.45'
$ws.Cells.Item(26, 4).Style = $savedStyle
$ws.Cells.Item(26, 5).Value = '  -0.17%  '

$ws.Cells.Item(27, 5).Value = '  -0.52%  '

$ws.Cells.Item(28, 5).Value = '  +0.17%  '

$savedStyle = $ws.Cells.Item(29, 4).Style
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '6.32'
$ws.Cells.Item(29, 4).Style = $savedStyle
$ws.Cells.Item(29, 5).Value = '  -1.91%  '

$ws.Cells.Item(30, 5).Value = '  -0.14%  '

$ws.Cells.Item(31, 5).Value = '  +2.16%  '

$ws.Cells.Item(32, 5).Value = '  -4.05%  '

$ws.Cells.Item(33, 5).Value = '  -0.69%  '

$ws.Cells.Item(34, 5).Value = '  -0.20%  '

$savedStyle = $ws.Cells.Item(35, 4).Style
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.387.08'
$ws.Cells.Item(35, 4).Style = $savedStyle
$ws.Cells.Item(35, 5).Value = '  -0.21%  '

$ws.Cells.Item(36, 5).Value = '  +2.19%  '

$ws.Cells.Item(37, 5).Value = '  -3.11%  '

$ws.Cells.Item(38, 5).Value = '  -0.31%  '

$ws.Cells.Item(39, 5).Value = '  +2.41%  '

$savedStyle = $ws.Cells.Item(40, 4).Style
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.0162'
$ws.Cells.Item(40, 4).Style = $savedStyle
$ws.Cells.Item(40, 5).Value = '  -1.72%  '

$ws.Cells.Item(41, 2).Value = 'RenderToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$savedStyle = $ws.Cells.Item(41, 4).Style
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '1.95'
$ws.Cells.Item(41, 4).Style = $savedStyle
$ws.Cells.Item(41, 5).Value = '  +3.19%  '

$ws.Cells.Item(42, 2).Value = 'ImmutableX'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$savedStyle = $ws.Cells.Item(42, 4).Style
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.520'
$ws.Cells.Item(42, 4).Style = $savedStyle
$ws.Cells.Item(42, 5).Value = '  -2.99%  '

$ws.Cells.Item(43, 5).Value = '  -0.12%  '

$ws.Cells.Item(44, 5).Value = '  -0.37%  '

$ws.Cells.Item(45, 5).Value = '  -2.19%  '

$ws.Cells.Item(46, 5).Value = '  -3.18%  '

$savedStyle = $ws.Cells.Item(47, 4).Style
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '62.25'
$ws.Cells.Item(47, 4).Style = $savedStyle
$ws.Cells.Item(47, 5).Value = '  -0.73%  '

$savedStyle = $ws.Cells.Item(48, 4).Style
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.917'
$ws.Cells.Item(48, 4).Style = $savedStyle
$ws.Cells.Item(48, 5).Value = '  -6.07%  '

$savedStyle = $ws.Cells.Item(49, 4).Style
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '1.701.66'
$ws.Cells.Item(49, 4).Style = $savedStyle
$ws.Cells.Item(49, 5).Value = '  -0.06%  '

$savedStyle = $ws.Cells.Item(50, 4).Style
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '85.44'
$ws.Cells.Item(50, 4).Style = $savedStyle
$ws.Cells.Item(50, 5).Value = '  -0.69%  '

$ws.Cells.Item(51, 5).Value = '  -1.93%  '
